$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 11365158
$ws.Range("I15").Value = 11365158
$ws.Range("K15").Value = 34095474
$ws.Range("M15").Value = -34095305
$ws.Range("H112").Value = 1335.3
$ws.Range("J112").Value = 1367.6
$ws.Range("L112").Value = 4102.799999999999
$ws.Range("N112").Value = -6318.799999999999
$ws.Range("H113").Value = 3126.2104
$ws.Range("J113").Value = 2475.25
$ws.Range("L113").Value = 2475.25
$ws.Range("N113").Value = -8983.25
$ws.Range("H116").Value = 3602.875
$ws.Range("I116").Value = 3649.0908
$ws.Range("J116").Value = 3501.2
$ws.Range("K116").Value = 3649.0908
$ws.Range("L116").Value = 3501.2
$ws.Range("M116").Value = -207.0907999999999
$ws.Range("N116").Value = -10385.2
$ws.Range("H137").Value = 969
$ws.Range("I137").Value = 817
$ws.Range("J137").Value = 1982.3334
$ws.Range("K137").Value = 2451
$ws.Range("L137").Value = 5947.0002
$ws.Range("M137").Value = 99
$ws.Range("N137").Value = -11047.0002
$ws.Range("H138").Value = 2726.78
$ws.Range("I138").Value = 632.8276
$ws.Range("J138").Value = 3582.0564
$ws.Range("K138").Value = 1898.4828
$ws.Range("L138").Value = 10746.1692
$ws.Range("M138").Value = 3241.5172
$ws.Range("N138").Value = -21026.1692
$ws.Range("H139").Value = 50524.285
$ws.Range("J139").Value = 50524.285
$ws.Range("L139").Value = 50524.285
$ws.Range("N139").Value = -60804.285
$ws.Range("H141").Value = 3720.3572
$ws.Range("I141").Value = 2716.3635
$ws.Range("K141").Value = 8149.0905
$ws.Range("M141").Value = -2969.0905

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 7131.476
$ws.Range("I32").Value = 6441
$ws.Range("K32").Value = 6441
$ws.Range("M32").Value = -6154
$ws.Range("H44").Value = 22049
$ws.Range("J44").Value = 22049
$ws.Range("L44").Value = 22049
$ws.Range("N44").Value = -23025
$ws.Range("H55").Value = 21056.4
$ws.Range("J55").Value = 21056.4
$ws.Range("L55").Value = 21056.4
$ws.Range("N55").Value = -21686.4
$ws.Range("H61").Value = 1540.375
$ws.Range("I61").Value = 1418.64
$ws.Range("J61").Value = 1975.1428
$ws.Range("K61").Value = 1418.64
$ws.Range("L61").Value = 1975.1428
$ws.Range("M61").Value = -1206.64
$ws.Range("N61").Value = -2399.1428
$ws.Range("H63").Value = 1999
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("N63").ClearContents()
$ws.Range("H66").Value = 1999
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("N66").ClearContents()
$ws.Range("H80").Value = 21389.092
$ws.Range("J80").Value = 21389.092
$ws.Range("L80").Value = 21389.092
$ws.Range("N80").Value = -23385.092
$ws.Range("H83").Value = 21389.092
$ws.Range("J83").Value = 21389.092
$ws.Range("L83").Value = 64167.276
$ws.Range("N83").Value = -74151.276
$ws.Range("H126").Value = 7710
$ws.Range("I126").Value = 7710
$ws.Range("K126").Value = 23130
$ws.Range("M126").Value = -20660
$ws.Range("H136").Value = 1540.375
$ws.Range("I136").Value = 1418.64
$ws.Range("J136").Value = 1975.1428
$ws.Range("K136").Value = 4255.92
$ws.Range("L136").Value = 5925.428400000001
$ws.Range("M136").Value = -1705.92
$ws.Range("N136").Value = -11025.4284

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H19").Value = 19990
$ws.Range("I19").Value = 19990
$ws.Range("K19").Value = 19990
$ws.Range("M19").Value = -19817
$ws.Range("H35").Value = 26900
$ws.Range("H82").Value = 14637.167
$ws.Range("I82").Value = 8128.5
$ws.Range("J82").Value = 17891.5
$ws.Range("K82").Value = 8128.5
$ws.Range("L82").Value = 17891.5
$ws.Range("M82").Value = -7745.5
$ws.Range("N82").Value = -18657.5
$ws.Range("H85").Value = 14637.167
$ws.Range("I85").Value = 8128.5
$ws.Range("J85").Value = 17891.5
$ws.Range("K85").Value = 8128.5
$ws.Range("L85").Value = 17891.5
$ws.Range("M85").Value = -6802.5
$ws.Range("N85").Value = -20543.5
$ws.Range("H94").Value = 980.4286
$ws.Range("I94").Value = 652.1667
$ws.Range("J94").Value = 2950
$ws.Range("K94").Value = 652.1667
$ws.Range("L94").Value = 2950
$ws.Range("M94").Value = -201.1667
$ws.Range("N94").Value = -3852
$ws.Range("H128").Value = 2461.25
$ws.Range("I128").Value = 2461.25
$ws.Range("K128").Value = 7383.75
$ws.Range("M128").Value = -4893.75
$ws.Range("H134").Value = 26011.592
$ws.Range("I134").Value = 3598.375
$ws.Range("K134").Value = 10795.125
$ws.Range("M134").Value = -8260.125

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H4").Value = 0
$ws.Range("J4").Value = 0
$ws.Range("L4").Value = 0
$ws.Range("N4").ClearContents()
$ws.Range("H16").Value = 973.9
$ws.Range("I16").Value = 948.7778
$ws.Range("J16").Value = 1200
$ws.Range("K16").Value = 948.7778
$ws.Range("L16").Value = 1200
$ws.Range("M16").Value = -661.7778
$ws.Range("N16").Value = -1774
$ws.Range("H31").Value = 3005.2856
$ws.Range("I31").Value = 3069
$ws.Range("J31").Value = 2400
$ws.Range("K31").Value = 3069
$ws.Range("L31").Value = 2400
$ws.Range("M31").Value = -2774
$ws.Range("N31").Value = -2990
$ws.Range("H34").Value = 3005.2856
$ws.Range("I34").Value = 3069
$ws.Range("J34").Value = 2400
$ws.Range("K34").Value = 3069
$ws.Range("L34").Value = 2400
$ws.Range("M34").Value = -2867
$ws.Range("N34").Value = -2804
$ws.Range("H41").Value = 6300.857
$ws.Range("I41").Value = 2821.6
$ws.Range("J41").Value = 14999
$ws.Range("K41").Value = 2821.6
$ws.Range("L41").Value = 14999
$ws.Range("M41").Value = -2393.6
$ws.Range("N41").Value = -15855
$ws.Range("H113").Value = 973.9
$ws.Range("I113").Value = 948.7778
$ws.Range("J113").Value = 1200
$ws.Range("K113").Value = 948.7778
$ws.Range("L113").Value = 1200
$ws.Range("M113").Value = 1221.2222
$ws.Range("N113").Value = -5540

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 664.3200000000001
$ws.Range("I5").Value = 654.26086
$ws.Range("K5").Value = 1962.78258
$ws.Range("M5").Value = -1850.78258
$ws.Range("H9").Value = 51191524
$ws.Range("I9").Value = 86685180
$ws.Range("J9").Value = 33444694
$ws.Range("K9").Value = 260055540
$ws.Range("L9").Value = 100334082
$ws.Range("M9").Value = -260055316
$ws.Range("N9").Value = -100334530
$ws.Range("H40").Value = 105
$ws.Range("I40").Value = 105
$ws.Range("K40").Value = 420
$ws.Range("M40").Value = -351
$ws.Range("H50").Value = 86.666664
$ws.Range("I50").Value = 86.666664
$ws.Range("K50").Value = 259.999992
$ws.Range("M50").Value = 221.000008
$ws.Range("H53").Value = 86.666664
$ws.Range("I53").Value = 86.666664
$ws.Range("K53").Value = 259.999992
$ws.Range("M53").Value = 221.000008
$ws.Range("H121").Value = 613.1
$ws.Range("I121").Value = 439.75
$ws.Range("J121").Value = 656.4375
$ws.Range("K121").Value = 1319.25
$ws.Range("L121").Value = 1969.3125
$ws.Range("M121").Value = -9.25
$ws.Range("N121").Value = -4589.3125
$ws.Range("H122").Value = 345683.9
$ws.Range("I122").Value = 548.375
$ws.Range("J122").Value = 770466.0600000001
$ws.Range("K122").Value = 4935.375
$ws.Range("L122").Value = 6934194.540000001
$ws.Range("M122").Value = -2485.375
$ws.Range("N122").Value = -6939094.540000001
$ws.Range("H131").Value = 6037626
$ws.Range("I131").Value = 41751224
$ws.Range("J131").Value = 1525.1831
$ws.Range("K131").Value = 125253672
$ws.Range("L131").Value = 4575.5493
$ws.Range("M131").Value = -125248632
$ws.Range("N131").Value = -14655.5493
$ws.Range("H132").Value = 701.6667
$ws.Range("I132").Value = 750
$ws.Range("J132").Value = 605
$ws.Range("K132").Value = 6750
$ws.Range("L132").Value = 5445
$ws.Range("M132").Value = -4220
$ws.Range("N132").Value = -10505
$ws.Range("H135").Value = 664.3200000000001
$ws.Range("I135").Value = 654.26086
$ws.Range("K135").Value = 5888.34774
$ws.Range("M135").Value = -3353.34774

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H20").Value = 70006
$ws.Range("J20").Value = 70006
$ws.Range("L20").Value = 70006
$ws.Range("N20").Value = -70496
$ws.Range("H126").Value = 3409.8823
$ws.Range("I126").Value = 3453.3333
$ws.Range("J126").Value = 3305.6
$ws.Range("K126").Value = 10359.9999
$ws.Range("L126").Value = 9916.799999999999
$ws.Range("M126").Value = -7889.999899999999
$ws.Range("N126").Value = -14856.8

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 1760.5
$ws.Range("I40").Value = 1808.0435
$ws.Range("J40").Value = 1604.2858
$ws.Range("K40").Value = 1808.0435
$ws.Range("L40").Value = 1604.2858
$ws.Range("M40").Value = -1672.0435
$ws.Range("N40").Value = -1876.2858
$ws.Range("H43").Value = 502500
$ws.Range("I43").Value = 1000000
$ws.Range("J43").Value = 5000
$ws.Range("K43").Value = 1000000
$ws.Range("L43").Value = 5000
$ws.Range("M43").Value = -999807
$ws.Range("N43").Value = -5386

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H18").Value = 2000
$ws.Range("J18").Value = 2000
$ws.Range("L18").Value = 2000
$ws.Range("N18").Value = -2346
$ws.Range("H20").Value = 1399.75
$ws.Range("I20").Value = 999.5
$ws.Range("J20").Value = 1800
$ws.Range("K20").Value = 999.5
$ws.Range("L20").Value = 1800
$ws.Range("M20").Value = -759.5
$ws.Range("N20").Value = -2280
$ws.Range("H107").Value = 583.80646
$ws.Range("I107").Value = 657.3913
$ws.Range("J107").Value = 372.25
$ws.Range("K107").Value = 1972.1739
$ws.Range("L107").Value = 1116.75
$ws.Range("M107").Value = -52.1739
$ws.Range("N107").Value = -4956.75
$ws.Range("H123").Value = 28967.53
$ws.Range("J123").Value = 28967.53
$ws.Range("L123").Value = 28967.53
$ws.Range("N123").Value = -38767.53
